$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.124.05"
$ws.Range("E2").Value = "  +0.74%  "
$ws.Range("D3").Value = "2.267.32"
$ws.Range("E3").Value = "  -0.31%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "230.99"
$ws.Range("E5").Value = "  -0.12%  "
$ws.Range("D6").Value = "0.639"
$ws.Range("E6").Value = "  +1.71%  "
$ws.Range("D7").Value = "64.02"
$ws.Range("E7").Value = "  +4.20%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").Value = "0.447"
$ws.Range("E9").Value = "  +5.63%  "
$ws.Range("E10").Value = "  +5.90%  "
$ws.Range("E11").Value = "  -0.72%  "
$ws.Range("D12").Value = "27.31"
$ws.Range("E13").Value = "  +1.77%  "
$ws.Range("D14").Value = "2.605.55"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "15.79"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("E16").Value = "  +5.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.840"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = "  +3.56%  "
$ws.Range("D18").Value = "2.272.55"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "43.896.94"
$ws.Range("E19").Value = "  +0.29%  "
$ws.Range("E20").Value = "  +7.37%  "
$ws.Range("D21").Value = "73.88"
$ws.Range("E21").Value = "  +1.01%  "
$ws.Range("E22").Value = "  -1.96%  "
$ws.Range("D23").Value = "252.95"
$ws.Range("E23").Value = "  +0.43%  "
$ws.Range("E24").Value = "  +0.10%  "
$ws.Range("E25").Value = "  -4.18%  "
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D26").Value = "10.12"
$ws.Range("E26").Value = "  +2.48%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "2.26"
$ws.Range("E27").Value = "  -4.70%  "
$ws.Range("D28").Value = "3.32"
$ws.Range("E28").Value = "  +24.45%  "
$ws.Range("D29").Value = "171.41"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.140"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -0.89%  "
$ws.Range("D31").Value = "20.91"
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("E32").Value = "  -4.30%  "
$ws.Range("E33").Value = "  +2.72%  "
$ws.Range("D34").Value = "0.0708"
$ws.Range("E34").Value = "  +6.88%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").Value = "4.89"
$ws.Range("E36").Value = "  -3.43%  "
$ws.Range("E37").Value = "  +5.60%  "
$ws.Range("D38").Value = "6.52"
$ws.Range("E38").Value = "  +0.44%  "
$ws.Range("D39").Value = "2.32"
$ws.Range("E39").Value = "  -3.69%  "
$ws.Range("E40").Value = "  +3.33%  "
$ws.Range("E41").Value = "  +0.16%  "
$ws.Range("E42").Value = "  +0.33%  "
$ws.Range("D43").Value = "0.0988"
$ws.Range("E43").Value = "  -0.06%  "
$ws.Range("D44").Value = "17.49"
$ws.Range("E44").Value = "  +5.12%  "
$ws.Range("E45").Value = "  -5.93%  "
$ws.Range("D46").Value = "10.43"
$ws.Range("E46").Value = "  +8.80%  "
$ws.Range("E47").Value = "  -0.46%  "
$ws.Range("D48").Value = "98.45"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("D49").Value = "4.38"
$ws.Range("E49").Value = "  -3.36%  "
$ws.Range("D50").Value = "2.34"
$ws.Range("E50").Value = "  +2.83%  "
$ws.Range("D51").Value = "1.446.99"
$ws.Range("E51").Value = "  -1.60%  "
